# Note added and DB fixed
# Adds a new "Lesser Hydra" creature row to the stat-dice source sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17

$ws.Cells.Item($row, 1).Value  = "Lesser Hydra"   # A17 Creature
$ws.Cells.Item($row, 2).Value  = "2D6+12"         # B17 STRSpec
$ws.Cells.Item($row, 3).Value  = "1D6+12"         # C17 CONSpec
$ws.Cells.Item($row, 4).Value  = "2D6+18"         # D17 SIZSpec
$ws.Cells.Item($row, 5).Value  = "2D6+12"         # E17 DEXSpec
$ws.Cells.Item($row, 7).Value  = "2D6+6"          # G17 POWSpec
$ws.Cells.Item($row, 15).Value = 6                # O17 Move
$ws.Cells.Item($row, 21).Value = "Lesser Hydra"   # U17 Hit_location
